$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook receives its weekly update: two new "Tuna" price records
# (most recent week) are inserted right before the current row 24, pushing
# every existing record down by two rows (old 24-39 -> new 26-41).
$ws.Rows("24:25").Insert()

# New row 24: Especial quality record
$ws.Cells.Item(24, 1).Value = 5
$ws.Cells.Item(24, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(24, 3).Value = "Maule"
$ws.Cells.Item(24, 4).Value = 44985
$ws.Cells.Item(24, 5).Value = 7
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100107
$ws.Cells.Item(24, 8).Value = "Otros"
$ws.Cells.Item(24, 9).Value = 100107011
$ws.Cells.Item(24, 10).Value = "Tuna"
$ws.Cells.Item(24, 11).Value = "Sin especificar"
$ws.Cells.Item(24, 12).Value = "Especial"
$ws.Cells.Item(24, 13).Value = 300
$ws.Cells.Item(24, 14).Value = 18000
$ws.Cells.Item(24, 15).Value = 18000
$ws.Cells.Item(24, 16).Value = 18000
$ws.Cells.Item(24, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(24, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(24, 19).Value = 1000
$ws.Cells.Item(24, 20).Value = 18

# New row 25: Segunda quality record
$ws.Cells.Item(25, 1).Value = 5
$ws.Cells.Item(25, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(25, 3).Value = "Maule"
$ws.Cells.Item(25, 4).Value = 44985
$ws.Cells.Item(25, 5).Value = 7
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100107
$ws.Cells.Item(25, 8).Value = "Otros"
$ws.Cells.Item(25, 9).Value = 100107011
$ws.Cells.Item(25, 10).Value = "Tuna"
$ws.Cells.Item(25, 11).Value = "Sin especificar"
$ws.Cells.Item(25, 12).Value = "Segunda"
$ws.Cells.Item(25, 13).Value = 150
$ws.Cells.Item(25, 14).Value = 12000
$ws.Cells.Item(25, 15).Value = 12000
$ws.Cells.Item(25, 16).Value = 12000
$ws.Cells.Item(25, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(25, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(25, 19).Value = 667
$ws.Cells.Item(25, 20).Value = 18

# Make sure the date column keeps the same date/time number format used by
# the rest of the column (Insert normally carries it down already, but set
# it explicitly to be safe).
$ws.Range("D24:D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
